# Commit: "pres and time investigation"
#
# On the "inference" sheet, the "w stim and isi time" row (row 41) had its
# per-run padding multiplier changed from 32*4 to 32*2.5. B41 holds the
# formula directly; C41:E41 share it (fill-right from C41). Updating the
# formulas here lets the workbook's own dependency chain (rows 42, 43, 45,
# 46, 47, 48) recompute to the new values automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inference")

$ws.Range("B41").Formula = "=B40+32*2.5"
$ws.Range("C41:E41").Formula = "=C40+32*2.5"

# Match the author's final scroll position / selection on that sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
[void]$ws.Range("I45").Select()
